$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 91 ("「美女と野獣」" entry), shifting all subsequent rows up by one.
$ws.Rows.Item(91).Delete()
